$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.060111333333333
$ws.Range("H2").Value = 6.180334
$ws.Range("I2").Value = 0.2095457297481522
$ws.Range("J2").Value = 0.2095457297481522
$ws.Range("M2").Value = 0.8077876666666667
$ws.Range("N2").Value = 2.423363
$ws.Range("O2").Value = 0.04902038147436601
$ws.Range("P2").Value = 0.04902038147436601
$ws.Range("Q2").Value = 1.664132527026889
$ws.Range("R2").Value = 14.977192743242
$ws.Range("S2").Value = 0.01027201160857882
$ws.Range("T2").Value = 0.01027201160857883
$ws.Range("G3").Value = 2.060111333333333
$ws.Range("H3").Value = 6.180334
$ws.Range("I3").Value = 0.2095457297481522
$ws.Range("J3").Value = 0.2095457297481522
$ws.Range("O3").Value = 0.4722854529078861
$ws.Range("P3").Value = 0.4722854529078861
$ws.Range("Q3").Value = 16.03303688357933
$ws.Range("R3").Value = 144.297331952214
$ws.Range("S3").Value = 0.09896539987901955
$ws.Range("T3").Value = 0.09896539987901956
$ws.Range("G4").Value = 2.060111333333333
$ws.Range("H4").Value = 6.180334
$ws.Range("I4").Value = 0.2095457297481522
$ws.Range("J4").Value = 0.2095457297481522
$ws.Range("M4").Value = 5.009781333333333
$ws.Range("N4").Value = 15.029344
$ws.Range("O4").Value = 0.3040172587389813
$ws.Range("P4").Value = 0.3040172587389813
$ws.Range("Q4").Value = 10.32070730232178
$ws.Range("R4").Value = 92.88636572089599
$ws.Range("S4").Value = 0.06370551833849263
$ws.Range("T4").Value = 0.06370551833849263
$ws.Range("G5").Value = 2.060111333333333
$ws.Range("H5").Value = 6.180334
$ws.Range("I5").Value = 0.2095457297481522
$ws.Range("J5").Value = 0.2095457297481522
$ws.Range("M5").Value = 2.878432333333333
$ws.Range("N5").Value = 8.635297
$ws.Range("O5").Value = 0.1746769068787666
$ws.Range("P5").Value = 0.1746769068787666
$ws.Range("Q5").Value = 5.929891072133111
$ws.Range("R5").Value = 53.369019649198
$ws.Range("S5").Value = 0.03660279992206115
$ws.Range("T5").Value = 0.03660279992206116
$ws.Range("I6").Value = 0.431709024657012
$ws.Range("J6").Value = 0.431709024657012
$ws.Range("M6").Value = 0.8077876666666667
$ws.Range("N6").Value = 2.423363
$ws.Range("O6").Value = 0.04902038147436601
$ws.Range("P6").Value = 0.04902038147436601
$ws.Range("Q6").Value = 3.428468960003334
$ws.Range("R6").Value = 30.85622064003
$ws.Range("S6").Value = 0.02116254107461321
$ws.Range("T6").Value = 0.02116254107461321
$ws.Range("I7").Value = 0.431709024657012
$ws.Range("J7").Value = 0.431709024657012
$ws.Range("O7").Value = 0.4722854529078861
$ws.Range("P7").Value = 0.4722854529078861
$ws.Range("S7").Value = 0.2038898922345587
$ws.Range("T7").Value = 0.2038898922345587
$ws.Range("I8").Value = 0.431709024657012
$ws.Range("J8").Value = 0.431709024657012
$ws.Range("M8").Value = 5.009781333333333
$ws.Range("N8").Value = 15.029344
$ws.Range("O8").Value = 0.3040172587389813
$ws.Range("P8").Value = 0.3040172587389813
$ws.Range("Q8").Value = 21.26286461962667
$ws.Range("R8").Value = 191.36578157664
$ws.Range("S8").Value = 0.1312469942491041
$ws.Range("T8").Value = 0.1312469942491041
$ws.Range("I9").Value = 0.431709024657012
$ws.Range("J9").Value = 0.431709024657012
$ws.Range("M9").Value = 2.878432333333333
$ws.Range("N9").Value = 8.635297
$ws.Range("O9").Value = 0.1746769068787666
$ws.Range("P9").Value = 0.1746769068787666
$ws.Range("Q9").Value = 12.21684399939667
$ws.Range("R9").Value = 109.95159599457
$ws.Range("S9").Value = 0.07540959709873601
$ws.Range("T9").Value = 0.07540959709873601
$ws.Range("G10").Value = 2.602283
$ws.Range("H10").Value = 7.806849
$ws.Range("I10").Value = 0.2646931170287289
$ws.Range("J10").Value = 0.2646931170287289
$ws.Range("M10").Value = 0.8077876666666667
$ws.Range("N10").Value = 2.423363
$ws.Range("O10").Value = 0.04902038147436601
$ws.Range("P10").Value = 0.04902038147436601
$ws.Range("Q10").Value = 2.102092112576333
$ws.Range("R10").Value = 18.918829013187
$ws.Range("S10").Value = 0.01297535757038729
$ws.Range("T10").Value = 0.0129753575703873
$ws.Range("G11").Value = 2.602283
$ws.Range("H11").Value = 7.806849
$ws.Range("I11").Value = 0.2646931170287289
$ws.Range("J11").Value = 0.2646931170287289
$ws.Range("O11").Value = 0.4722854529078861
$ws.Range("P11").Value = 0.4722854529078861
$ws.Range("Q11").Value = 20.252545891781
$ws.Range("R11").Value = 182.272913026029
$ws.Range("S11").Value = 0.1250107086575133
$ws.Range("T11").Value = 0.1250107086575133
$ws.Range("G12").Value = 2.602283
$ws.Range("H12").Value = 7.806849
$ws.Range("I12").Value = 0.2646931170287289
$ws.Range("J12").Value = 0.2646931170287289
$ws.Range("M12").Value = 5.009781333333333
$ws.Range("N12").Value = 15.029344
$ws.Range("O12").Value = 0.3040172587389813
$ws.Range("P12").Value = 0.3040172587389813
$ws.Range("Q12").Value = 13.03686879745066
$ws.Range("R12").Value = 117.331819177056
$ws.Range("S12").Value = 0.08047127584615052
$ws.Range("T12").Value = 0.08047127584615052
$ws.Range("G13").Value = 2.602283
$ws.Range("H13").Value = 7.806849
$ws.Range("I13").Value = 0.2646931170287289
$ws.Range("J13").Value = 0.2646931170287289
$ws.Range("M13").Value = 2.878432333333333
$ws.Range("N13").Value = 8.635297
$ws.Range("O13").Value = 0.1746769068787666
$ws.Range("P13").Value = 0.1746769068787666
$ws.Range("Q13").Value = 7.490495527683666
$ws.Range("R13").Value = 67.414459749153
$ws.Range("S13").Value = 0.04623577495467773
$ws.Range("T13").Value = 0.04623577495467773
$ws.Range("G14").Value = 0.9246566666666666
$ws.Range("H14").Value = 2.77397
$ws.Range("I14").Value = 0.09405212856610688
$ws.Range("J14").Value = 0.0940521285661069
$ws.Range("M14").Value = 0.8077876666666667
$ws.Range("N14").Value = 2.423363
$ws.Range("O14").Value = 0.04902038147436601
$ws.Range("P14").Value = 0.04902038147436601
$ws.Range("Q14").Value = 0.7469262512344443
$ws.Range("R14").Value = 6.72233626111
$ws.Range("S14").Value = 0.004610471220786676
$ws.Range("T14").Value = 0.004610471220786677
$ws.Range("G15").Value = 0.9246566666666666
$ws.Range("H15").Value = 2.77397
$ws.Range("I15").Value = 0.09405212856610688
$ws.Range("J15").Value = 0.0940521285661069
$ws.Range("O15").Value = 0.4722854529078861
$ws.Range("P15").Value = 0.4722854529078861
$ws.Range("Q15").Value = 7.196239446596666
$ws.Range("R15").Value = 64.76615501936999
$ws.Range("S15").Value = 0.04441945213679452
$ws.Range("T15").Value = 0.04441945213679453
$ws.Range("G16").Value = 0.9246566666666666
$ws.Range("H16").Value = 2.77397
$ws.Range("I16").Value = 0.09405212856610688
$ws.Range("J16").Value = 0.0940521285661069
$ws.Range("M16").Value = 5.009781333333333
$ws.Range("N16").Value = 15.029344
$ws.Range("O16").Value = 0.3040172587389813
$ws.Range("P16").Value = 0.3040172587389813
$ws.Range("Q16").Value = 4.632327708408888
$ws.Range("R16").Value = 41.69094937567999
$ws.Range("S16").Value = 0.02859347030523405
$ws.Range("T16").Value = 0.02859347030523405
$ws.Range("G17").Value = 0.9246566666666666
$ws.Range("H17").Value = 2.77397
$ws.Range("I17").Value = 0.09405212856610688
$ws.Range("J17").Value = 0.0940521285661069
$ws.Range("M17").Value = 2.878432333333333
$ws.Range("N17").Value = 8.635297
$ws.Range("O17").Value = 0.1746769068787666
$ws.Range("P17").Value = 0.1746769068787666
$ws.Range("Q17").Value = 2.661561646565555
$ws.Range("R17").Value = 23.95405481909
$ws.Range("S17").Value = 0.01642873490329163
$ws.Range("T17").Value = 0.01642873490329163
